# Add files via upload
# - Populate ID/MSHS/Email columns (C, D, E) for rows 2-28 on sheet "MUL10L3"
# - Reset selection to C2 on all other sheets
# - Widen column F on sheet "9H02" to fit long names

$wb = $excel.ActiveWorkbook

# --- 1) Reset cursor/selection back to C2 on the other six sheets ---
$otherSheets = @("12T1", "12T2", "9H01", "9H02", "9H03", "9H07")
foreach ($name in $otherSheets) {
    $s = $wb.Worksheets.Item($name)
    $s.Activate() | Out-Null
    $s.Range("C2").Select() | Out-Null
}

# --- 2) Widen column F (Ho va ten) on 9H02 so long names are fully visible ---
$ws9h02 = $wb.Worksheets.Item("9H02")
$ws9h02.Columns.Item(6).ColumnWidth = 34.66

# --- 3) Fill in Ma dinh danh (C), MSHS (D) and Email (E) for the MUL10L3 roster ---
$ws = $wb.Worksheets.Item("MUL10L3")

$data = @(
    @{ Row=2; C="7939735860"; D=2107018; E="2107018@lsts.edu.vn" },
    @{ Row=3; C="7982708662"; D=2110001; E="2110001@lsts.edu.vn" },
    @{ Row=4; C="7926735798"; D=2106077; E="2106077@lsts.edu.vn" },
    @{ Row=5; C="7925828318"; D=2106043; E="2106043@lsts.edu.vn" },
    @{ Row=6; C="7965951207"; D=2106157; E="2106157@lsts.edu.vn" },
    @{ Row=7; C="7949028877"; D=2110016; E="2110016@lsts.edu.vn" },
    @{ Row=8; C="7965951096"; D=2106050; E="2106050@lsts.edu.vn" },
    @{ Row=9; C="7941660363"; D=2107019; E="2107019@lsts.edu.vn" },
    @{ Row=10; C="7965951098"; D=2106221; E="2106221@lsts.edu.vn" },
    @{ Row=11; C="7965951059"; D=2106160; E="2106160@lsts.edu.vn" },
    @{ Row=12; C="7937595064"; D=2108015; E="2108015@lsts.edu.vn" },
    @{ Row=13; C="7924539916"; D=2106158; E="2106158@lsts.edu.vn" },
    @{ Row=14; C="7966247118"; D=2107012; E="2107012@lsts.edu.vn" },
    @{ Row=15; C="7926730716"; D=2107020; E="2107020@lsts.edu.vn" },
    @{ Row=16; C="7965951256"; D=2106203; E="2106203@lsts.edu.vn" },
    @{ Row=17; C="7927755049"; D=2106019; E="2106019@lsts.edu.vn" },
    @{ Row=18; C="7925401212"; D=2107033; E="2107033@lsts.edu.vn" },
    @{ Row=19; C="7962071478"; D=2107003; E="2107003@lsts.edu.vn" },
    @{ Row=20; C="7963199074"; D=2108002; E="2108002@lsts.edu.vn" },
    @{ Row=21; C="3465951222"; D=2106171; E="2106171@lsts.edu.vn" },
    @{ Row=22; C="7926731173"; D=2106095; E="2106095@lsts.edu.vn" },
    @{ Row=23; C="7928261265"; D=2106205; E="2106205@lsts.edu.vn" },
    @{ Row=24; C="7988530875"; D=2110030; E="2110030@lsts.edu.vn" },
    @{ Row=25; C="7951465925"; D=2106219; E="2106219@lsts.edu.vn" },
    @{ Row=26; C="7942595748"; D=2110032; E="2110032@lsts.edu.vn" },
    @{ Row=27; C="7965951151"; D=2106099; E="2106099@lsts.edu.vn" },
    @{ Row=28; C="8006607113"; D=2107049; E="2107049@lsts.edu.vn" }
)

$firstRow = $data[0].Row
$lastRow = $data[$data.Count - 1].Row

# Mark column C as Text first so 10-digit phone numbers are stored as strings
# (not silently coerced to numbers), matching the source roster formatting.
$ws.Range("C$firstRow`:C$lastRow").NumberFormat = "@"

foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 4).Value = $item.D
}
foreach ($item in $data) {
    $ws.Cells.Item($item.Row, 5).Value = $item.E
}

# Strip the temporary Text number format back off so the cells keep the
# workbook's default (unstyled) look, same as every other populated column.
$ws.Range("C$firstRow`:C$lastRow").ClearFormats()

# --- 4) Leave the MUL10L3 tab as the active sheet/selection (C2), as before ---
$ws.Activate() | Out-Null
$ws.Range("C2").Select() | Out-Null

